$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC23_Verify_UserRegistration")

# Insert a new row at position 4 (a WAIT step before the CLICK RegisterButton step),
# pushing all existing rows down by one
$ws.Rows("4:4").Insert()

# Give the new row the same thin-bordered look as the rest of the data table
$ws.Range("A4:E4").Borders.LineStyle = 1

# Populate the newly inserted row with the WAIT keyword in column B
$ws.Range("B4").Value = "WAIT"

# Mirror the selection Excel leaves behind after this kind of edit
$ws.Range("C3").Select()
